$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-7 from 45243 (2023-11-13)
# to 45244 (2023-11-14), keeping existing date formatting.
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 3).Value = 45244
}
